$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("J2").Value = 1.18
$ws.Range("K2").Value = 4.5
$ws.Range("L2").Value = 1.83
$ws.Range("M2").Value = 1.83
$ws.Range("N2").Value = 4
$ws.Range("O2").Value = 1.25
$ws.Range("P2").Value = 1.85
$ws.Range("Q2").Value = 1.95

# Row 3
$ws.Range("G3").Value = 2.05
$ws.Range("H3").Value = 3.4
$ws.Range("I3").Value = 3.7
$ws.Range("U3").Value = 8.5
$ws.Range("X3").Value = 19
$ws.Range("AE3").Value = 8.5

# Row 4
$ws.Range("I4").Value = 2.75
$ws.Range("J4").Value = 1.1
$ws.Range("K4").Value = 7
$ws.Range("L4").Value = 1.44
$ws.Range("M4").Value = 2.63
$ws.Range("N4").Value = 2.35
$ws.Range("O4").Value = 1.57
$ws.Range("P4").Value = 1.53
$ws.Range("Q4").Value = 2.38
$ws.Range("R4").Value = 2
$ws.Range("S4").Value = 1.73
$ws.Range("T4").Value = 7
$ws.Range("Y4").Value = 41
$ws.Range("Z4").Value = 7
$ws.Range("AC4").Value = 67
$ws.Range("AF4").Value = 12
$ws.Range("AI4").Value = 26

# Row 5
$ws.Range("G5").Value = 2.9
$ws.Range("J5").Value = 1.03
$ws.Range("K5").Value = 17
$ws.Range("P5").Value = 1.29
$ws.Range("Q5").Value = 3.5
$ws.Range("W5").Value = 34
$ws.Range("Y5").Value = 23
$ws.Range("AA5").Value = 7.5

# Row 6
$ws.Range("G6").Value = 4.75
$ws.Range("H6").Value = 3.8
$ws.Range("I6").Value = 1.67
$ws.Range("N6").Value = 1.8
$ws.Range("O6").Value = 2
$ws.Range("R6").Value = 1.75
$ws.Range("S6").Value = 2
$ws.Range("U6").Value = 26
$ws.Range("X6").Value = 41
$ws.Range("Y6").Value = 41
$ws.Range("AB6").Value = 15
$ws.Range("AD6").Value = 201
$ws.Range("AF6").Value = 8.5
$ws.Range("AH6").Value = 13

# Row 7
$ws.Range("G7").Value = 1.75
$ws.Range("H7").Value = 4
$ws.Range("I7").Value = 4.33
$ws.Range("W7").Value = 15
$ws.Range("Z7").Value = 13
$ws.Range("AA7").Value = 7.5
$ws.Range("AE7").Value = 13
$ws.Range("AF7").Value = 23
$ws.Range("AG7").Value = 13
$ws.Range("AH7").Value = 41

# Row 9
$ws.Range("I9").Value = 3.25

# Row 10
$ws.Range("J10").Value = 1.08
$ws.Range("K10").Value = 8

# Row 11
$ws.Range("G11").Value = 1.85
$ws.Range("H11").Value = 3.2
$ws.Range("I11").Value = 4
$ws.Range("J11").Value = 1.08
$ws.Range("K11").Value = 8
$ws.Range("R11").Value = 2
$ws.Range("S11").Value = 1.73
$ws.Range("T11").Value = 6
$ws.Range("U11").Value = 8
$ws.Range("W11").Value = 15
$ws.Range("AC11").Value = 67
$ws.Range("AF11").Value = 21
$ws.Range("AI11").Value = 41

# Row 13
$ws.Range("G13").Value = 3.6
$ws.Range("I13").Value = 1.95
$ws.Range("K13").Value = 7.5
$ws.Range("AG13").Value = 9.5

# Row 16
$ws.Range("N16").Value = 1.44
$ws.Range("R16").Value = 2.75
$ws.Range("S16").Value = 1.4

# Row 18
$ws.Range("G18").Value = 2.1
$ws.Range("H18").Value = 3.3
$ws.Range("J18").Value = 1.06
$ws.Range("K18").Value = 10
$ws.Range("L18").Value = 1.29
$ws.Range("M18").Value = 3.5
$ws.Range("N18").Value = 2
$ws.Range("O18").Value = 1.85
$ws.Range("P18").Value = 1.4
$ws.Range("Q18").Value = 2.75
$ws.Range("R18").Value = 1.75
$ws.Range("S18").Value = 2
$ws.Range("T18").Value = 8
$ws.Range("U18").Value = 10
$ws.Range("X18").Value = 17
$ws.Range("Y18").Value = 26
$ws.Range("Z18").Value = 10
$ws.Range("AD18").Value = 201
$ws.Range("AE18").Value = 11
$ws.Range("AI18").Value = 29

# Row 19
$ws.Range("G19").Value = 1.53
$ws.Range("H19").Value = 4.33
$ws.Range("L19").Value = 1.29
$ws.Range("M19").Value = 3.5
$ws.Range("N19").Value = 1.95
$ws.Range("O19").Value = 1.9
$ws.Range("R19").Value = 2
$ws.Range("S19").Value = 1.75
$ws.Range("T19").Value = 6.5
$ws.Range("U19").Value = 7
$ws.Range("Z19").Value = 11
$ws.Range("AD19").Value = 401
$ws.Range("AE19").Value = 13

# Row 20
$ws.Range("N20").Value = 2.15
$ws.Range("O20").Value = 1.67
